$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for Price column (D) to avoid Excel coercing numeric-looking
# strings (e.g. "3.70", "0.996") into floating point numbers, which would strip
# meaningful trailing zeros / change the stored cell type from text to number.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.067.70"
$ws.Range("E2").Value = "  +2.75%  "
$ws.Range("D3").Value = "1.581.95"
$ws.Range("E3").Value = "  +1.89%  "
$ws.Range("D4").Value = "0.996"
$ws.Range("E4").Value = "  -0.60%  "
$ws.Range("D5").Value = "211.92"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").Value = "0.518"
$ws.Range("E6").Value = "  +7.01%  "
$ws.Range("D7").Value = "0.995"
$ws.Range("E7").Value = "  -0.63%  "
$ws.Range("D8").Value = "25.49"
$ws.Range("E8").Value = "  +8.55%  "
$ws.Range("E9").Value = "  +2.76%  "
$ws.Range("D10").Value = "0.0594"
$ws.Range("E10").Value = "  +1.82%  "
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").Value = "1.809.13"
$ws.Range("E12").Value = "  +2.00%  "
$ws.Range("D13").Value = "1.611.64"
$ws.Range("E13").Value = "  +3.77%  "
$ws.Range("D14").Value = "29.092.37"
$ws.Range("E14").Value = "  +2.80%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "3.70"
$ws.Range("E15").Value = "  +2.08%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "0.521"
$ws.Range("E16").Value = "  +2.50%  "
$ws.Range("D17").Value = "62.45"
$ws.Range("E17").Value = "  +3.26%  "
$ws.Range("D18").Value = "238.49"
$ws.Range("E18").Value = "  +5.17%  "
$ws.Range("D19").Value = "7.42"
$ws.Range("E19").Value = "  +1.61%  "
$ws.Range("D20").Value = "0.0₃0694"
$ws.Range("E20").Value = "  +2.87%  "
$ws.Range("D21").Value = "0.994"
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("E22").Value = "  +2.23%  "
$ws.Range("D23").Value = "9.18"
$ws.Range("E23").Value = "  +4.14%  "
$ws.Range("D24").Value = "2.10"
$ws.Range("E24").Value = "  +4.95%  "
$ws.Range("D25").Value = "152.93"
$ws.Range("E25").Value = "  +3.42%  "
$ws.Range("E26").Value = "  +4.91%  "
$ws.Range("D27").Value = "15.14"
$ws.Range("E27").Value = "  +2.58%  "
$ws.Range("D28").Value = "6.32"
$ws.Range("E28").Value = "  +1.64%  "
$ws.Range("D29").Value = "0.996"
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("E30").Value = "  -0.68%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("E32").Value = "  +1.41%  "
$ws.Range("D33").Value = "1.421.79"
$ws.Range("E33").Value = "  +2.65%  "
$ws.Range("D34").Value = "3.04"
$ws.Range("E34").Value = "  -0.57%  "
$ws.Range("E35").Value = "  -1.49%  "
$ws.Range("E36").Value = "  +0.98%  "
$ws.Range("E37").Value = "  +7.78%  "
$ws.Range("D38").Value = "2.28"
$ws.Range("E38").Value = "  -2.03%  "
$ws.Range("E39").Value = "  +1.09%  "
$ws.Range("E40").Value = "  +2.67%  "
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("B42").Value = "BitcoinSV"
$ws.Range("C42").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D42").Value = "53.46"
$ws.Range("E42").Value = "  +26.83%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "0.995"
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "0.787"
$ws.Range("E44").Value = "  +1.34%  "
$ws.Range("E45").Value = "  -1.42%  "
$ws.Range("D46").Value = "64.64"
$ws.Range("E46").Value = "  +4.61%  "
$ws.Range("D47").Value = "5.33"
$ws.Range("E47").Value = "  -1.68%  "
$ws.Range("D48").Value = "1.719.79"
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("D49").Value = "0.848"
$ws.Range("E49").Value = "  -6.43%  "
$ws.Range("D50").Value = "85.51"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("E51").Value = "  +0.72%  "

# Clean up: remove the temporary Text number format so the cells end up with
# no explicit style, matching the original (unstyled) Price column cells.
$priceRange.ClearFormats()
